$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values scraped for this run,
# as reflected by the commit "Updated cryptos list ... with GitHub Actions".
# A $null value means that column is unchanged for that row.
$updates = @(
    @{ Row = 2; D = "30.446.38"; E = "  -0.93%  " },
    @{ Row = 3; D = "1.891.51"; E = "  -1.02%  " },
    @{ Row = 4; D = "1.000"; E = "  -0.23%  " },
    @{ Row = 5; D = "240.20"; E = "  +0.29%  " },
    @{ Row = 6; D = "0.9999"; E = "  -0.26%  " },
    @{ Row = 7; D = "0.4847"; E = "  -1.53%  " },
    @{ Row = 8; D = "0.2894"; E = "  -2.27%  " },
    @{ Row = 9; D = "0.06633"; E = "  -1.59%  " },
    @{ Row = 10; D = "1.899.62"; E = "  -0.76%  " },
    @{ Row = 11; D = "16.99"; E = "  -0.37%  " },
    @{ Row = 12; D = "0.07411"; E = "  +0.62%  " },
    @{ Row = 13; D = "5.206"; E = "  +0.88%  " },
    @{ Row = 14; D = "89.31"; E = $null },
    @{ Row = 15; D = "0.6652"; E = "  -0.74%  " },
    @{ Row = 16; D = "30.418.35"; E = "  -0.89%  " },
    @{ Row = 17; D = "13.56"; E = "  +0.68%  " },
    @{ Row = 18; D = "0.000007792"; E = "  -1.36%  " },
    @{ Row = 19; D = "0.9993"; E = "  -0.25%  " },
    @{ Row = 20; D = "5.437"; E = "  +2.46%  " },
    @{ Row = 21; D = "2.143.71"; E = "  -1.18%  " },
    @{ Row = 22; D = "1.000"; E = "  -0.22%  " },
    @{ Row = 23; D = "218.08"; E = "  +11.63%  " },
    @{ Row = 24; D = "6.212"; E = $null },
    @{ Row = 25; D = "9.454"; E = "  -1.84%  " },
    @{ Row = 26; D = "165.10"; E = "  +1.01%  " },
    @{ Row = 27; D = "18.64"; E = "  +0.27%  " },
    @{ Row = 28; D = "1.949"; E = "  +0.06%  " },
    @{ Row = 29; D = $null; E = "  -2.41%  " },
    @{ Row = 30; D = "4.326"; E = "  -1.21%  " },
    @{ Row = 31; D = "0.09198"; E = "  +0.95%  " },
    @{ Row = 32; D = $null; E = "  +0.51%  " },
    @{ Row = 33; D = $null; E = "  -3.10%  " },
    @{ Row = 34; D = "0.7509"; E = "  +1.65%  " },
    @{ Row = 35; D = "1.159"; E = "  +4.29%  " },
    @{ Row = 36; D = "2.704"; E = "  -0.76%  " },
    @{ Row = 37; D = "0.01904"; E = "  +4.66%  " },
    @{ Row = 38; D = "2.642"; E = "  -2.71%  " },
    @{ Row = 39; D = "0.9211"; E = "  +0.00%  " },
    @{ Row = 40; D = "2.099"; E = "  +1.26%  " },
    @{ Row = 41; D = "6.092"; E = "  +3.19%  " },
    @{ Row = 42; D = "107.55"; E = "  +0.67%  " },
    @{ Row = 43; D = "0.4362"; E = "  -1.73%  " },
    @{ Row = 44; D = "1.002"; E = "  +0.13%  " },
    @{ Row = 45; D = "7.661"; E = "  +0.95%  " },
    @{ Row = 46; D = "0.1350"; E = "  -2.37%  " },
    @{ Row = 47; D = "66.20"; E = "  -11.51%  " },
    @{ Row = 48; D = "1.580"; E = "  +10.57%  " },
    @{ Row = 49; D = "8.961"; E = "  -1.09%  " },
    @{ Row = 50; D = "34.46"; E = "  -2.96%  " },
    @{ Row = 51; D = "0.05701"; E = "  -2.69%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Price strings (e.g. "1.000", "30.446.38") look numeric/date-like to
        # Excel's automatic type conversion, so force the cell to Text first,
        # assign the literal string, then restore the default "Normal" style
        # so the cell's style index matches the untouched workbook cells.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        # Volume strings are already safe text (padded with spaces, contain
        # "%"), so they can be assigned directly without losing formatting.
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
